# estimates_PS2025.xlsx: refresh the rolling "infVol" (column I) series.
# Three new trailing observations were appended upstream to the source data,
# so the existing infVol history shifts down by three rows (I2:I255 -> I5:I258)
# and the three newly opened rows at the top (I2:I4) get freshly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current infVol column (rows 2-255) before overwriting anything,
# then shift it down into rows 5-258.
$srcRange = $ws.Range("I2:I255")
$srcValues = $srcRange.Value2

$dstRange = $ws.Range("I5:I258")
$dstRange.Value2 = $srcValues

# Newly computed leading values for the three rows this freed up at the top.
$ws.Range("I2").Value2 = 0.720864933637685
$ws.Range("I3").Value2 = 0.69123318704604797
$ws.Range("I4").Value2 = 0.67916760993880199

# Reflect the active selection on the frozen bottom-right pane.
$ws.Range("F1").Select()
